# Updates cryptos list values (Price / Volume(1h)) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.702.44"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.190.99"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.54"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.97"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.188.06"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.509"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.27"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "3.712.05"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "65.834.35"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.30"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "3.194.78"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.23"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.81"
$ws.Range("E21").Value = "  +5.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.17"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.96"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.38"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.19"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.86"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.99"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.68"
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.44"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0916"
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "485.30"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0423"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.89"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "3.012.25"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.289"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").Value = "0.0₃0639"
$ws.Range("E46").Value = "  +7.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.06"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.60"
$ws.Range("E51").Value = "  -2.37%  "
